$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold/centered/bordered) onto the new I1:J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF)
$data = @(
    @(2, 1, 3),
    @(3, 1, 5),
    @(4, 1, 4),
    @(5, 1, 5),
    @(6, 1, 5),
    @(7, 1, 4),
    @(8, 1, 4),
    @(9, 1, 6),
    @(10, 7, 7),
    @(11, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
